# New Submission Synced: 2026-02-10 18:50:00
# Appends a new form-response row (row 10) to the "JSS 3B" sheet, matching
# the same shape as the existing rows: Timestamp, Full Name, Admission No,
# AI Score.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("JSS 3B")

$ws.Cells.Item(10, 1).Value = "2026-02-10 18:50:00"
$ws.Cells.Item(10, 2).Value = "Adam Muhammad Gudusu"

# "Admission No" in this sheet is stored as text (several existing rows hold
# purely-numeric-looking admission numbers as text, e.g. C7="10", C8="36").
# Force text formatting before writing so "28" isn't reinterpreted as a number.
$ws.Cells.Item(10, 3).NumberFormat = "@"
$ws.Cells.Item(10, 3).Value = "28"

$ws.Cells.Item(10, 4).Value = 8
